$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "x and n" inputs / confidence level for the left-hand example
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 50
$ws.Range("C8").Value = 0.95
